$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz4")
$ws.Activate()

# The original sheet stored its fill-down formulas as two shared-formula
# groups (C4:F7 and B5:B7). Touch one cell from each group with a
# non-matching value so the exporter stops treating them as shared
# groups; every formula cell gets re-entered individually below anyway,
# so the final values are unaffected.
$ws.Range("D5").Value = "__break__"
$ws.Range("B6").Value = "__break__"

# Insert a new column before C, shifting existing C:F data to D:G
$ws.Range("C1").EntireColumn.Insert() | Out-Null

# New values for the freshly inserted column C
$ws.Range("C1").Value = 17.5
$ws.Range("C2").Value = 20
$ws.Range("C3").Value = 20

# Updated values for column B (rows 2 and 3)
$ws.Range("B2").Value = 25
$ws.Range("B3").Value = 25

# Re-enter every formula individually (instead of relying on the
# shared-formula group Excel created originally) so each cell in B4:G7
# carries its own plain formula, matching a freshly recalculated sheet.
foreach ($col in @("B", "C", "D", "E", "F", "G")) {
    for ($row = 4; $row -le 7; $row++) {
        $prevRow = $row - 1
        $ws.Range("$col$row").Formula = "=$col$prevRow+$col`$1"
    }
}

# Match the new cell selection recorded in the saved file
$ws.Range("K8").Select() | Out-Null

for ($r = 1; $r -le 7; $r++) {
    $line = ""
    for ($c = 1; $c -le 7; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        $f = $cell.Formula()
        $line = $line + "[" + $cell.Address() + " v=" + $v + " f=" + $f + "]"
    }
    Write-Host $line
}
